# Updates the "cryptos" price table (Coin/Link/Price/Volume(1h)) with the
# latest scraped values. Source data stores numeric-looking prices (e.g.
# "229.11") as plain text, so a leading apostrophe is used on values that
# Excel would otherwise auto-convert to a Number when assigned via .Value;
# this keeps those cells as text, matching the original data format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '38.132.53'
$ws.Cells.Item(2, 5).Value = '  +0.46%  '

$ws.Cells.Item(3, 4).Value = '2.043.16'
$ws.Cells.Item(3, 5).Value = '  -0.56%  '

$ws.Cells.Item(4, 5).Value = '  +0.10%  '

$ws.Cells.Item(5, 4).Value = '''229.11'
$ws.Cells.Item(5, 5).Value = '  -0.26%  '

$ws.Cells.Item(6, 4).Value = '''0.615'
$ws.Cells.Item(6, 5).Value = '  -0.36%  '

$ws.Cells.Item(7, 4).Value = '''61.05'
$ws.Cells.Item(7, 5).Value = '  +3.99%  '

$ws.Cells.Item(8, 5).Value = '  +0.01%  '

$ws.Cells.Item(9, 4).Value = '''0.384'
$ws.Cells.Item(9, 5).Value = '  -0.48%  '

$ws.Cells.Item(10, 4).Value = '''0.0822'
$ws.Cells.Item(10, 5).Value = '  +1.53%  '

$ws.Cells.Item(11, 5).Value = '  +0.51%  '

$ws.Cells.Item(12, 4).Value = '''14.73'
$ws.Cells.Item(12, 5).Value = '  +0.35%  '

$ws.Cells.Item(13, 4).Value = '2.345.28'
$ws.Cells.Item(13, 5).Value = '  -0.50%  '

$ws.Cells.Item(14, 4).Value = '''21.44'
$ws.Cells.Item(14, 5).Value = '  +2.67%  '

$ws.Cells.Item(15, 4).Value = '''0.771'
$ws.Cells.Item(15, 5).Value = '  +2.35%  '

$ws.Cells.Item(16, 5).Value = '  -1.46%  '

$ws.Cells.Item(17, 4).Value = '2.047.37'
$ws.Cells.Item(17, 5).Value = '  +0.25%  '

$ws.Cells.Item(18, 4).Value = '38.042.14'
$ws.Cells.Item(18, 5).Value = '  +0.37%  '

$ws.Cells.Item(19, 4).Value = '''6.00'
$ws.Cells.Item(19, 5).Value = '  -4.53%  '

$ws.Cells.Item(20, 4).Value = '''70.10'
$ws.Cells.Item(20, 5).Value = '  +0.54%  '

$ws.Cells.Item(21, 5).Value = '  -1.04%  '

$ws.Cells.Item(22, 4).Value = '''225.16'
$ws.Cells.Item(22, 5).Value = '  +0.34%  '

$ws.Cells.Item(23, 4).Value = '''1.00'
$ws.Cells.Item(23, 5).Value = '  +0.04%  '

$ws.Cells.Item(24, 5).Value = '  +0.22%  '

$ws.Cells.Item(25, 4).Value = '''2.26'
$ws.Cells.Item(25, 5).Value = '  +0.13%  '

$ws.Cells.Item(26, 2).Value = 'Cosmos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(26, 4).Value = '''9.35'
$ws.Cells.Item(26, 5).Value = '  +0.52%  '

$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(27, 4).Value = '''167.59'
$ws.Cells.Item(27, 5).Value = '  +0.75%  '

$ws.Cells.Item(28, 5).Value = '  -2.79%  '

$ws.Cells.Item(29, 4).Value = '''19.00'
$ws.Cells.Item(29, 5).Value = '  -0.18%  '

$ws.Cells.Item(30, 5).Value = '  -2.67%  '

$ws.Cells.Item(31, 4).Value = '''0.121'
$ws.Cells.Item(31, 5).Value = '  +1.04%  '

$ws.Cells.Item(32, 5).Value = '  +5.28%  '

$ws.Cells.Item(33, 4).Value = '''4.45'
$ws.Cells.Item(33, 5).Value = '  -1.90%  '

$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Value = '''0.0611'
$ws.Cells.Item(34, 5).Value = '  +0.13%  '

$ws.Cells.Item(35, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(35, 4).Value = '''4.56'
$ws.Cells.Item(35, 5).Value = '  -0.96%  '

$ws.Cells.Item(36, 5).Value = '  +6.24%  '

$ws.Cells.Item(37, 5).Value = '  -1.67%  '

$ws.Cells.Item(38, 4).Value = '''3.30'
$ws.Cells.Item(38, 5).Value = '  +0.26%  '

$ws.Cells.Item(39, 5).Value = '  -0.20%  '

$ws.Cells.Item(40, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(40, 4).Value = '''17.88'
$ws.Cells.Item(40, 5).Value = '  +7.94%  '

$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(41, 4).Value = '1.527.17'
$ws.Cells.Item(41, 5).Value = '  +2.01%  '

$ws.Cells.Item(42, 4).Value = '''0.0219'
$ws.Cells.Item(42, 5).Value = '  +0.78%  '

$ws.Cells.Item(43, 4).Value = '''96.81'
$ws.Cells.Item(43, 5).Value = '  -0.49%  '

$ws.Cells.Item(44, 4).Value = '''2.82'
$ws.Cells.Item(44, 5).Value = '  -2.20%  '

$ws.Cells.Item(45, 4).Value = '''0.0921'
$ws.Cells.Item(45, 5).Value = '  +0.22%  '

$ws.Cells.Item(46, 4).Value = '''1.11'
$ws.Cells.Item(46, 5).Value = '  -2.60%  '

$ws.Cells.Item(47, 4).Value = '''4.04'
$ws.Cells.Item(47, 5).Value = '  -1.72%  '

$ws.Cells.Item(48, 4).Value = '''1.02'
$ws.Cells.Item(48, 5).Value = '  -0.08%  '

$ws.Cells.Item(49, 5).Value = '  +0.15%  '

$ws.Cells.Item(50, 4).Value = '''7.15'
$ws.Cells.Item(50, 5).Value = '  +0.68%  '

$ws.Cells.Item(51, 4).Value = '2.233.73'
$ws.Cells.Item(51, 5).Value = '  -0.57%  '
